# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.789.39'
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").Value = '2.044.00'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '227.47'
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("E6").Value = '  -0.65%  '

$ws.Range("D7").Value = '60.02'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '0.377'
$ws.Range("E9").Value = '  -2.36%  '

$ws.Range("D10").Value = '0.0844'
$ws.Range("E10").Value = '  +3.06%  '

$ws.Range("E11").Value = '  -0.03%  '

$ws.Range("D12").Value = '2.346.93'
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("E13").Value = '  -1.45%  '

$ws.Range("D14").Value = '21.10'
$ws.Range("E14").Value = '  +0.34%  '

$ws.Range("D15").Value = '5.50'
$ws.Range("E15").Value = '  +5.78%  '

$ws.Range("E16").Value = '  +0.76%  '

$ws.Range("D17").Value = '2.037.58'
$ws.Range("E17").Value = '  -0.24%  '

$ws.Range("D18").Value = '37.788.09'
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("D19").Value = '69.45'
$ws.Range("E19").Value = '  -0.47%  '

$ws.Range("D20").Value = '5.93'
$ws.Range("E20").Value = '  -2.26%  '

$ws.Range("D21").Value = '0.0₃0826'
$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").Value = '223.32'
$ws.Range("E22").Value = '  -0.76%  '

$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("E25").Value = '  +3.11%  '

$ws.Range("D26").Value = '169.65'
$ws.Range("E26").Value = '  +2.73%  '

$ws.Range("E27").Value = '  +1.05%  '

$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("D29").Value = '18.78'
$ws.Range("E29").Value = '  -0.84%  '

$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("E31").Value = '  -0.64%  '

$ws.Range("D32").Value = '2.25'
$ws.Range("E32").Value = '  +8.73%  '

$ws.Range("D33").Value = '4.38'
$ws.Range("E33").Value = '  -1.28%  '

$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").Value = '6.57'
$ws.Range("E36").Value = '  +2.50%  '

$ws.Range("E37").Value = '  +3.80%  '

$ws.Range("D38").Value = '3.44'
$ws.Range("E38").Value = '  +6.06%  '

$ws.Range("E39").Value = '  -0.10%  '

$ws.Range("D40").Value = '17.93'
$ws.Range("E40").Value = '  +6.46%  '

$ws.Range("D41").Value = '1.528.70'
$ws.Range("E41").Value = '  -0.95%  '

$ws.Range("D42").Value = '97.88'
$ws.Range("E42").Value = '  +0.92%  '

$ws.Range("E43").Value = '  -0.89%  '

$ws.Range("E44").Value = '  -0.06%  '

$ws.Range("E45").Value = '  -1.91%  '

$ws.Range("D46").Value = '4.17'
$ws.Range("E46").Value = '  +6.48%  '

$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("D49").Value = '2.94'
$ws.Range("E49").Value = '  -0.47%  '

$ws.Range("D50").Value = '7.09'
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("D51").Value = '2.236.49'
$ws.Range("E51").Value = '  +0.45%  '
